$wb = $excel.ActiveWorkbook

# Cell value updates per sheet, derived from the authoritative diff of the
# committed OOXML (scheduled market-price refresh for the Alexander Profits sheet).


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 27784726
$ws.Range("I62").Value = 1619.6154
$ws.Range("J62").Value = 100020800
$ws.Range("K62").Value = 1619.6154
$ws.Range("L62").Value = 100020800
$ws.Range("M62").Value = -995.6153999999999
$ws.Range("N62").Value = -100022048
$ws.Range("H65").Value = 27784726
$ws.Range("I65").Value = 1619.6154
$ws.Range("J65").Value = 100020800
$ws.Range("K65").Value = 8098.076999999999
$ws.Range("L65").Value = 500104000
$ws.Range("M65").Value = -4978.076999999999
$ws.Range("N65").Value = -500110240
$ws.Range("H118").Value = 84001730
$ws.Range("J118").Value = 2796
$ws.Range("L118").Value = 8388
$ws.Range("N118").Value = -11702
$ws.Range("H129").Value = 703.0833
$ws.Range("I129").Value = 432.09525
$ws.Range("J129").Value = 2600
$ws.Range("K129").Value = 1296.28575
$ws.Range("L129").Value = 7800
$ws.Range("M129").Value = 3703.71425
$ws.Range("N129").Value = -17800

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1805.4897
$ws.Range("I2").Value = 1443.0294
$ws.Range("J2").Value = 2627.0667
$ws.Range("K2").Value = 1443.0294
$ws.Range("L2").Value = 2627.0667
$ws.Range("M2").Value = -1330.0294
$ws.Range("N2").Value = -2853.0667
$ws.Range("H116").Value = 1805.4897
$ws.Range("I116").Value = 1443.0294
$ws.Range("J116").Value = 2627.0667
$ws.Range("K116").Value = 1443.0294
$ws.Range("L116").Value = 2627.0667
$ws.Range("M116").Value = 850.9706000000001
$ws.Range("N116").Value = -7215.066699999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1805.4897
$ws.Range("I3").Value = 1443.0294
$ws.Range("J3").Value = 2627.0667
$ws.Range("K3").Value = 1443.0294
$ws.Range("L3").Value = 2627.0667
$ws.Range("M3").Value = -1329.0294
$ws.Range("N3").Value = -2855.0667

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 251902.34
$ws.Range("I5").Value = 273.46667
$ws.Range("K5").Value = 820.4000100000001
$ws.Range("M5").Value = -708.4000100000001
$ws.Range("H22").Value = 605.5
$ws.Range("I22").Value = 469.23077
$ws.Range("J22").Value = 959.8
$ws.Range("K22").Value = 1407.69231
$ws.Range("L22").Value = 2879.4
$ws.Range("M22").Value = -1238.69231
$ws.Range("N22").Value = -3217.4
$ws.Range("H27").Value = 605.5
$ws.Range("I27").Value = 469.23077
$ws.Range("J27").Value = 959.8
$ws.Range("K27").Value = 1407.69231
$ws.Range("L27").Value = 2879.4
$ws.Range("M27").Value = -1305.69231
$ws.Range("N27").Value = -3083.4
$ws.Range("H44").Value = 659.21155
$ws.Range("I44").Value = 275
$ws.Range("J44").Value = 691.2292
$ws.Range("K44").Value = 825
$ws.Range("L44").Value = 2073.6876
$ws.Range("M44").Value = -427
$ws.Range("N44").Value = -2869.6876
$ws.Range("H58").Value = 2981.3513
$ws.Range("I58").Value = 1847.25
$ws.Range("J58").Value = 3118.818
$ws.Range("K58").Value = 5541.75
$ws.Range("L58").Value = 9356.454000000002
$ws.Range("M58").Value = -5413.75
$ws.Range("N58").Value = -9612.454000000002
$ws.Range("H64").Value = 4894.5
$ws.Range("I64").Value = 772.75
$ws.Range("J64").Value = 6162.731
$ws.Range("K64").Value = 2318.25
$ws.Range("L64").Value = 18488.193
$ws.Range("M64").Value = -2048.25
$ws.Range("N64").Value = -19028.193
$ws.Range("H67").Value = 4894.5
$ws.Range("I67").Value = 772.75
$ws.Range("J67").Value = 6162.731
$ws.Range("K67").Value = 2318.25
$ws.Range("L67").Value = 18488.193
$ws.Range("M67").Value = -1382.25
$ws.Range("N67").Value = -20360.193
$ws.Range("H69").Value = 2656.1538
$ws.Range("I69").Value = 433.33334
$ws.Range("J69").Value = 4561.4287
$ws.Range("K69").Value = 1300.00002
$ws.Range("L69").Value = 13684.2861
$ws.Range("M69").Value = -489.0000199999999
$ws.Range("N69").Value = -15306.2861
$ws.Range("H72").Value = 2656.1538
$ws.Range("I72").Value = 433.33334
$ws.Range("J72").Value = 4561.4287
$ws.Range("K72").Value = 3900.00006
$ws.Range("L72").Value = 41052.85830000001
$ws.Range("M72").Value = 155.9999399999997
$ws.Range("N72").Value = -49164.85830000001
$ws.Range("H76").Value = 7421.7393
$ws.Range("I76").Value = 2975
$ws.Range("J76").Value = 8357.895
$ws.Range("K76").Value = 8925
$ws.Range("L76").Value = 25073.685
$ws.Range("M76").Value = -8542
$ws.Range("N76").Value = -25839.685
$ws.Range("H79").Value = 7421.7393
$ws.Range("I79").Value = 2975
$ws.Range("J79").Value = 8357.895
$ws.Range("K79").Value = 8925
$ws.Range("L79").Value = 25073.685
$ws.Range("M79").Value = -7599
$ws.Range("N79").Value = -27725.685
$ws.Range("H94").Value = 6406.2856
$ws.Range("J94").Value = 6129.846
$ws.Range("L94").Value = 18389.538
$ws.Range("N94").Value = -19741.538
$ws.Range("H106").Value = 3274.2917
$ws.Range("J106").Value = 3274.2917
$ws.Range("L106").Value = 9822.875100000001
$ws.Range("N106").Value = -11714.8751
$ws.Range("H107").Value = 585.7143
$ws.Range("I107").Value = 258.63635
$ws.Range("J107").Value = 1785
$ws.Range("K107").Value = 775.90905
$ws.Range("L107").Value = 5355
$ws.Range("M107").Value = 1144.09095
$ws.Range("N107").Value = -9195
$ws.Range("H121").Value = 431267.03
$ws.Range("I121").Value = 150
$ws.Range("J121").Value = 460999.25
$ws.Range("K121").Value = 450
$ws.Range("L121").Value = 1382997.75
$ws.Range("M121").Value = 860
$ws.Range("N121").Value = -1385617.75
$ws.Range("H122").Value = 52298.39
$ws.Range("I122").Value = 368
$ws.Range("J122").Value = 56398.156
$ws.Range("K122").Value = 3312
$ws.Range("L122").Value = 507583.404
$ws.Range("M122").Value = -862
$ws.Range("N122").Value = -512483.404
$ws.Range("H123").Value = 5560.8335
$ws.Range("I123").Value = 3010
$ws.Range("J123").Value = 6411.1113
$ws.Range("K123").Value = 9030
$ws.Range("L123").Value = 19233.3339
$ws.Range("M123").Value = -6580
$ws.Range("N123").Value = -24133.3339
$ws.Range("H131").Value = 968.3684
$ws.Range("I131").Value = 532.5
$ws.Range("J131").Value = 1008.4483
$ws.Range("K131").Value = 1597.5
$ws.Range("L131").Value = 3025.3449
$ws.Range("M131").Value = 3442.5
$ws.Range("N131").Value = -13105.3449
$ws.Range("H135").Value = 251902.34
$ws.Range("I135").Value = 273.46667
$ws.Range("K135").Value = 2461.20003
$ws.Range("M135").Value = 73.79997000000003

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1426.7391
$ws.Range("I136").Value = 1320.2439
$ws.Range("J136").Value = 2300
$ws.Range("K136").Value = 3960.7317
$ws.Range("L136").Value = 6900
$ws.Range("M136").Value = -1410.7317
$ws.Range("N136").Value = -12000

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1655.6923
$ws.Range("I107").Value = 1046.5625
$ws.Range("K107").Value = 3139.6875
$ws.Range("M107").Value = -1219.6875
$ws.Range("H113").Value = 29392.6
$ws.Range("I113").Value = 52980.527
$ws.Range("J113").Value = 1381.9375
$ws.Range("K113").Value = 158941.581
$ws.Range("L113").Value = 4145.8125
$ws.Range("M113").Value = -156771.581
$ws.Range("N113").Value = -8485.8125
$ws.Range("H136").Value = 1620.18
$ws.Range("I136").Value = 1651.8096
$ws.Range("J136").Value = 1454.125
$ws.Range("K136").Value = 4955.4288
$ws.Range("L136").Value = 4362.375
$ws.Range("M136").Value = -2405.4288
$ws.Range("N136").Value = -9462.375
